# Fruta / hortaliza, semanal
#
# A new weekly record is inserted at row 15 (pushing the existing rows
# 15-47 down to 16-48). The sheet's used range grows from A1:T47 to
# A1:T48 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting rows 15-47 down to 16-48.
$ws.Rows.Item(15).EntireRow.Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Femacal de La Calera"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44560
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100108
$ws.Range("H15").Value = "Tropicales y subtropicales"
$ws.Range("I15").Value = 100108004
$ws.Range("J15").Value = "Papaya"
$ws.Range("K15").Value = "Cultivar IV Región"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 56
$ws.Range("N15").Value = 24000
$ws.Range("O15").Value = 24000
$ws.Range("P15").Value = 24000
$ws.Range("Q15").Value = "$/bandeja 10 kilos"
$ws.Range("R15").Value = "Provincia del Elquí"
$ws.Range("S15").Value = 2400
$ws.Range("T15").Value = 10
